# Full overwrite with latest ShiftBookExcel2 files
# Populate the header row (row 1) of Sheet1 with the ShiftBookExcel2 column
# headings: timeStamp, action, studentID, name, date, shift, LIC, LIC verified

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$headers = @("timeStamp", "action", "studentID", "name", "date", "shift", "LIC", "LIC verified")

for ($i = 0; $i -lt $headers.Count; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# Match the author's saved cursor position (cell F10 selected) recorded in the
# sheet view of the edited workbook.
$ws.Range("F10").Select() | Out-Null
